$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("策略更新")

# Update status in E5 from "pending" to "passed"
$ws.Range("E5").Value = "passed"

# Update remark in F5: prefix each lot line with "(ok) "
$newRemark = "郑州日盘的so，如果今天晚上正常，明天部署`nBDY 2017/1/3 19:57:46`n日盘手数暂定如下：`n(ok) ta fl34 3`n(ok) zc fl34 1`n(ok) ma fl34 5`n(ok) sr fl36 1`n(ok) cf fl36 1`n(ok) ta fl36 1`n(ok) zc fl36 1`n(ok) cf fw10 1`n(ok) ta fw10 3`n(ok) zc fw10 2`n(ok) zc fd10 2`n(ok) ta fd10 3`n"
$ws.Range("F5").Value = $newRemark

# Update the active selection to C5 to match the saved view state
$ws.Range("C5").Select()
